$d = $word.ActiveDocument

function Get-ParaIndex($pattern) {
    $counter = 0
    foreach ($p in $d.Paragraphs) {
        $counter = $counter + 1
        if ($p.Range.Text -match $pattern) {
            return $counter
        }
    }
    return 0
}

# ---------------------------------------------------------------------------
# Hunk 1: insert a new bullet paragraph "Identify critical transitions in
# nanocosm systems and warning signals" right before the
# "Model communities with DEB model ..." paragraph. The word "nanocosm" is
# wrapped in spell-check proofErr markers (matching the style used elsewhere
# in the document) and the "_GoBack" bookmark (last-edit marker) ends up at
# the end of this new paragraph.
# ---------------------------------------------------------------------------
$modelIdx = Get-ParaIndex("^Model communities with DEB model parameterized")
$d.Paragraphs($modelIdx).Range.InsertParagraphBefore()

# InsertParagraphBefore turns the (former) "Model communities..." paragraph
# slot into the new, empty paragraph, and pushes "Model communities..." one
# slot further down - re-resolve by index after the mutation.
$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Identify critical transitions in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>nanocosm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> systems and warning signals</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs($modelIdx).Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# Hunk 2: the old "_GoBack" bookmark that used to sit in the middle of the
# "Coincidence of several small disturbances ..." paragraph (splitting it
# into two runs) is removed, and the two runs are merged back into a single
# run.
# ---------------------------------------------------------------------------
$coincidenceIdx = Get-ParaIndex("^Coincidence of several small disturbances")
$cr = $d.Paragraphs($coincidenceIdx).Range
$crNoMark = $d.Range($cr.Start, $cr.End - 1)
$coincidenceXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Coincidence of several small disturbances or few large disturbances are sufficient to move a system out of region of stability towards irreversible degradation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$crNoMark.InsertXML($coincidenceXml)

# ---------------------------------------------------------------------------
# Hunk 3: the "lastRenderedPageBreak" marker moves from the start of the
# "C, D extinct" run to the start of the "C vital, D extinct" run.
# ---------------------------------------------------------------------------
$vitalIdx = Get-ParaIndex("^C vital, D extinct")
$vr = $d.Paragraphs($vitalIdx).Range
$vrNoMark = $d.Range($vr.Start, $vr.End - 1)
$vitalXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>C vital, D extinct</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$vrNoMark.InsertXML($vitalXml)

$extinctIdx = Get-ParaIndex("^C, D extinct")
$er = $d.Paragraphs($extinctIdx).Range
$erNoMark = $d.Range($er.Start, $er.End - 1)
$extinctXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>C, D extinct</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$erNoMark.InsertXML($extinctXml)

Write-Output "edit complete"
